$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 619.8909
$ws.Range("J17").Value = 627.6667
$ws.Range("L17").Value = 1883.0001
$ws.Range("N17").Value = -2219.0001

$ws.Range("H40").Value = 125001980
$ws.Range("I40").Value = 1950
$ws.Range("J40").Value = 166668670
$ws.Range("K40").Value = 1950
$ws.Range("L40").Value = 166668670
$ws.Range("M40").Value = -1775
$ws.Range("N40").Value = -166669020

$ws.Range("H76").Value = 11277.23
$ws.Range("I76").Value = 13260
$ws.Range("J76").Value = 4668
$ws.Range("K76").Value = 13260
$ws.Range("L76").Value = 4668
$ws.Range("M76").Value = -12945
$ws.Range("N76").Value = -5298

$ws.Range("H79").Value = 11277.23
$ws.Range("I79").Value = 13260
$ws.Range("J79").Value = 4668
$ws.Range("K79").Value = 13260
$ws.Range("L79").Value = 4668
$ws.Range("M79").Value = -12168
$ws.Range("N79").Value = -6852

$ws.Range("H92").Value = 1121.1111
$ws.Range("I92").Value = 245.06667
$ws.Range("K92").Value = 245.06667
$ws.Range("M92").Value = 1002.93333

$ws.Range("H96").Value = 230.90909
$ws.Range("I96").Value = 201.375
$ws.Range("K96").Value = 604.125
$ws.Range("M96").Value = 768.875

$ws.Range("H129").Value = 1191.3
$ws.Range("I129").Value = 481.3
$ws.Range("J129").Value = 1427.9667
$ws.Range("K129").Value = 1443.9
$ws.Range("L129").Value = 4283.9001
$ws.Range("M129").Value = 3556.1
$ws.Range("N129").Value = -14283.9001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = 12
$ws.Range("N5").Value = -324

$ws.Range("H97").Value = 3403.0688
$ws.Range("I97").Value = 2479.4736
$ws.Range("J97").Value = 5157.9
$ws.Range("K97").Value = 2479.4736
$ws.Range("L97").Value = 5157.9
$ws.Range("M97").Value = -1983.4736
$ws.Range("N97").Value = -6149.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = -330

$ws.Range("H22").Value = 934.44446
$ws.Range("I22").Value = 917.1429000000001
$ws.Range("J22").Value = 995
$ws.Range("K22").Value = 917.1429000000001
$ws.Range("L22").Value = 995
$ws.Range("M22").Value = -744.1429000000001
$ws.Range("N22").Value = -1341

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 29.461538
$ws.Range("I7").Value = 23.454546
$ws.Range("J7").Value = 62.5
$ws.Range("K7").Value = 23.454546
$ws.Range("L7").Value = 62.5
$ws.Range("M7").Value = 89.54545400000001
$ws.Range("N7").Value = -288.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 741555.7
$ws.Range("J5").Value = 1667682.2
$ws.Range("L5").Value = 5003046.6
$ws.Range("N5").Value = -5003270.6

$ws.Range("H63").Value = 4100.0713
$ws.Range("I63").Value = 1800.2
$ws.Range("J63").Value = 5377.778
$ws.Range("K63").Value = 5400.6
$ws.Range("L63").Value = 16133.334
$ws.Range("M63").Value = -4651.6
$ws.Range("N63").Value = -17631.334

$ws.Range("H66").Value = 4100.0713
$ws.Range("I66").Value = 1800.2
$ws.Range("J66").Value = 5377.778
$ws.Range("K66").Value = 16201.8
$ws.Range("L66").Value = 48400.002
$ws.Range("M66").Value = -12457.8
$ws.Range("N66").Value = -55888.002

$ws.Range("H114").Value = 2082.1428
$ws.Range("I114").Value = 304.8
$ws.Range("J114").Value = 3069.5557
$ws.Range("K114").Value = 914.4000000000001
$ws.Range("L114").Value = 9208.667099999999
$ws.Range("M114").Value = 2339.6
$ws.Range("N114").Value = -15716.6671

$ws.Range("H117").Value = 6014.2607
$ws.Range("I117").Value = 870
$ws.Range("J117").Value = 7097.263
$ws.Range("K117").Value = 2610
$ws.Range("L117").Value = 21291.789
$ws.Range("M117").Value = 832
$ws.Range("N117").Value = -28175.789

$ws.Range("H121").Value = 476986.6
$ws.Range("I121").Value = 308.75
$ws.Range("J121").Value = 667657.75
$ws.Range("K121").Value = 926.25
$ws.Range("L121").Value = 2002973.25
$ws.Range("M121").Value = 383.75
$ws.Range("N121").Value = -2005593.25

$ws.Range("H135").Value = 741555.7
$ws.Range("J135").Value = 1667682.2
$ws.Range("L135").Value = 15009139.8
$ws.Range("N135").Value = -15014209.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2183.3333
$ws.Range("I2").Value = 1500
$ws.Range("K2").Value = 1500
$ws.Range("M2").Value = -1388

$ws.Range("H11").Value = 4950
$ws.Range("I11").Value = 4900
$ws.Range("J11").Value = 5000
$ws.Range("K11").Value = 4900
$ws.Range("L11").Value = 5000
$ws.Range("M11").Value = -4760
$ws.Range("N11").Value = -5280

$ws.Range("H13").Value = 4500
$ws.Range("I13").Value = 1000
$ws.Range("J13").Value = 8000
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 8000
$ws.Range("M13").Value = -860
$ws.Range("N13").Value = -8280

$ws.Range("H18").Value = 12790
$ws.Range("I18").Value = 9980
$ws.Range("J18").Value = 15600
$ws.Range("K18").Value = 9980
$ws.Range("L18").Value = 15600
$ws.Range("M18").Value = -9808
$ws.Range("N18").Value = -15944

$ws.Range("H22").Value = 892.5
$ws.Range("I22").Value = 561.5714
$ws.Range("J22").Value = 1057.9642
$ws.Range("K22").Value = 561.5714
$ws.Range("L22").Value = 1057.9642
$ws.Range("M22").Value = -266.5714
$ws.Range("N22").Value = -1647.9642

$ws.Range("H27").Value = 892.5
$ws.Range("I27").Value = 561.5714
$ws.Range("J27").Value = 1057.9642
$ws.Range("K27").Value = 561.5714
$ws.Range("L27").Value = 1057.9642
$ws.Range("M27").Value = -454.5714
$ws.Range("N27").Value = -1271.9642

$ws.Range("H46").Value = 1196.5264
$ws.Range("I46").Value = 1035.7693
$ws.Range("J46").Value = 1544.8334
$ws.Range("K46").Value = 1035.7693
$ws.Range("L46").Value = 1544.8334
$ws.Range("M46").Value = -847.7692999999999
$ws.Range("N46").Value = -1920.8334

$ws.Range("H55").Value = 1019.5
$ws.Range("I55").Value = 450
$ws.Range("J55").Value = 1304.25
$ws.Range("K55").Value = 450
$ws.Range("L55").Value = 1304.25
$ws.Range("M55").Value = -277
$ws.Range("N55").Value = -1650.25

$ws.Range("H133").Value = 47709.25
$ws.Range("J133").Value = 47709.25
$ws.Range("L133").Value = 47709.25
$ws.Range("N133").Value = -52769.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 30000
$ws.Range("J80").Value = 30000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -31996

$ws.Range("H81").Value = 9247.923000000001
$ws.Range("I81").Value = 25580.25
$ws.Range("J81").Value = 1989.1111
$ws.Range("K81").Value = 51160.5
$ws.Range("L81").Value = 3978.2222
$ws.Range("M81").Value = -50099.5
$ws.Range("N81").Value = -6100.2222

$ws.Range("H83").Value = 30000
$ws.Range("J83").Value = 30000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -99984

$ws.Range("H84").Value = 9247.923000000001
$ws.Range("I84").Value = 25580.25
$ws.Range("J84").Value = 1989.1111
$ws.Range("K84").Value = 255802.5
$ws.Range("L84").Value = 19891.111
$ws.Range("M84").Value = -250498.5
$ws.Range("N84").Value = -30499.111
